$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Cells changing from numeric to text ("0" / "***.*") need NumberFormat forced to text first ---
$ws.Range("C16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"

$ws.Range("C16").Value = "0"
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "***.*"
$ws.Range("C22").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"

# --- Restore original style (General/right-aligned text, style index 14) via format-only paste ---
$ws.Range("A14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)

# --- Cells changing from text to numeric ---
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100

$ws.Range("F15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# --- Remaining plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -58.620689655172
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 210
$ws.Range("K16").Value = -14.761904761904
$ws.Range("L16").Value = 17.763157894736
$ws.Range("M16").Value = 54.310344827586
$ws.Range("N16").Value = -84.394071490845
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 160
$ws.Range("K17").Value = 7.38255033557
$ws.Range("L17").Value = 16.788321167883
$ws.Range("M17").Value = 90.47619047619
$ws.Range("N17").Value = -37.5
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = -21.428571428571
$ws.Range("I18").Value = 222
$ws.Range("J18").Value = 227
$ws.Range("K18").Value = -2.202643171806
$ws.Range("L18").Value = 17.460317460317
$ws.Range("M18").Value = 7.766990291262
$ws.Range("N18").Value = -91.026677445432
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = -18.60465116279
$ws.Range("F19").Value = 116
$ws.Range("G19").Value = 152
$ws.Range("H19").Value = -23.684210526315
$ws.Range("I19").Value = 1441
$ws.Range("J19").Value = 1505
$ws.Range("K19").Value = -4.252491694352
$ws.Range("L19").Value = 47.794871794871
$ws.Range("M19").Value = 35.178236397748
$ws.Range("N19").Value = -54.369854338188
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 136
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = -9.333333333333
$ws.Range("L20").Value = -9.333333333333
$ws.Range("M20").Value = 76.623376623376
$ws.Range("N20").Value = -95.388267209223
$ws.Range("D21").Value = 54
$ws.Range("E21").Value = -25.925925925925
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 215
$ws.Range("H21").Value = -20
$ws.Range("I21").Value = 2153
$ws.Range("J21").Value = 2253
$ws.Range("K21").Value = -4.438526409232
$ws.Range("L21").Value = 33.147804576376
$ws.Range("M21").Value = 37.836107554417
$ws.Range("N21").Value = -78.508684368137
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -40
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = 66.666666666666
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("M23").Value = -4.166666666666
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 102
$ws.Range("E24").Value = -52.941176470588
$ws.Range("F24").Value = 200
$ws.Range("G24").Value = 346
$ws.Range("H24").Value = -42.196531791907
$ws.Range("I24").Value = 2671
$ws.Range("J24").Value = 3414
$ws.Range("K24").Value = -21.763327475102
$ws.Range("L24").Value = 29.283639883833
$ws.Range("M24").Value = 83.447802197802
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -3.225806451612
$ws.Range("I25").Value = 293
$ws.Range("J25").Value = 321
$ws.Range("K25").Value = -8.722741433021
$ws.Range("L25").Value = 0.342465753424
$ws.Range("M25").Value = -1.677852348993
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 23
$ws.Range("K26").Value = 64.285714285714
$ws.Range("L26").Value = 21.052631578947
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 76
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 8.571428571428
$ws.Range("L27").Value = -6.172839506172
$ws.Range("F30").Value = 3
$ws.Range("H30").Value = 50
$ws.Range("I30").Value = 14
$ws.Range("K30").Value = -48.148148148148
$ws.Range("L30").Value = 0
